# "aggiornamento fino a 13/03" - append 4 new daily rows (252-255) to the
# report sheet, continuing the existing date/count series through 13/03.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last existing data row is 251 (column A styled with style index 2:
# bold, centered, bordered, custom date/time number format). Copy that
# formatting down into the four new rows before writing the new values.
$ws.Range("A251").Copy()
$ws.Range("A252:A255").PasteSpecial(-4122)  # xlPasteFormats

$newRows = @(
    @(44326, 1, 1, 25.4323499491353),
    @(44327, 0, 1, 25.4323499491353),
    @(44328, 0, 1, 25.4323499491353),
    @(44329, 0, 1, 25.4323499491353)
)

$r = 252
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}
